# Add "None of the above" as a new option to the blood_type list sheet,
# and update the B-column data validation on the main sheet to match.

$wb = $excel.ActiveWorkbook

$wsList = $wb.Worksheets.Item("blood_type list")
$wsMain = $wb.Worksheets.Item("Export as TSV")

# Append the new list entry right after "O" (row 4 -> new row 5).
$wsList.Range("A5").Value = "None of the above"

# Refresh the dropdown validation on column B so it covers the new row
# and update the error message to mention the new option.
$validation = $wsMain.Range("B2:B1048576").Validation
$validation.Formula1 = "='blood_type list'!`$A`$1:`$A`$5"
$validation.ErrorMessage = "Value must be one of: A / B / AB / O / None of the above."
